$d = $word.ActiveDocument

# --- Add two new rows to the (only) table, mirroring the formatting of the
# --- existing last row (Rows.Add clones the last row's cell formatting). ---
$t = $d.Tables(1)

$row1 = $t.Rows.Add()
$row1.Cells(1).Range.Text = " 26/02/2022"
$row1.Cells(2).Range.Text = "4 Hours"
$row1.Cells(3).Range.Text = "Simulation – Objective 4"
$row1.Cells(4).Range.Text = "Added more properties to the opinions system, including new designations for fear and rivalry opinions between nations. Also added the ability for the AI to make decisions based on their relations to other nations and how certain actions would impact their standing."

$row2 = $t.Rows.Add()
$row2.Cells(1).Range.Text = "27/02/2022"
$row2.Cells(2).Range.Text = "2 Hours 20 Minutes"
$row2.Cells(3).Range.Text = "Simulation – Objective 4"
$row2.Cells(4).Range.Text = "Added the new “ally” property for opinions, as well as the modifiers for actions that would impact an ally. Also added more positive opinion changes."

# --- Insert four new blank paragraphs right after the table (before the
# --- existing run of blank paragraphs that precede "ADD MEETINGS"). We do
# --- this by locating the first blank paragraph right after the table and
# --- replacing its single paragraph-mark range with five bare <w:p/>
# --- paragraphs (net effect: +4 new blank paragraphs, landing exactly
# --- where the old one was). ---
$rng = $d.Content
$found = $rng.Find.Execute("ADD MEETINGS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$addMeetingsPara = $rng.Paragraphs(1)

$t = $d.Tables(1)
$n = $t.Rows.Count
Write-Host "Table now has $n rows"

$target = $addMeetingsPara.Previous()
$target = $target.Range.Duplicate
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$blankParaXml = "<w:p $ns/>"
$target.InsertXML($blankParaXml + $blankParaXml + $blankParaXml + $blankParaXml + $blankParaXml)

Write-Host "Paragraph count now: $($d.Paragraphs.Count)"
